$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: becomes the filled-in product data (was row 3 content), V2 updated to file name
$ws.Range("A2").Value = 'Summary of Product Characteristics'
$ws.Range("B2").Value = 'Awiqli'
$ws.Range("C2").Value = 'insulin icodec'
$ws.Range("D2").Value = '700 units/mL'
$ws.Range("E2").Value = 'Solution for injection in pre-filled pen'
$ws.Range("F2").Value = 'Subcutaneous use'
$ws.Range("G2").Value = 'Treatment of diabetes mellitus in adults'
$ws.Range("H2").Value = 'Once-weekly subcutaneous administration, dose range 10-700 units per injection'
$ws.Range("I2").Value = 'Hypersensitivity to the active substance or to any of the excipients'
$ws.Range("J2").Value = 'Risk of hypoglycaemia, hypersensitivity reactions, and immune system disorders'
$ws.Range("K2").Value = 'Novo Nordisk A/S'
$ws.Range("L2").Value = 'EU/1/24/1815/001, EU/1/24/1815/002, EU/1/24/1815/003, EU/1/24/1815/004, EU/1/24/1815/005, EU/1/24/1815/006, EU/1/24/1815/007, EU/1/24/1815/008, EU/1/24/1815/009, EU/1/24/1815/010, EU/1/24/1815/011, EU/1/24/1815/012, EU/1/24/1815/013, EU/1/24/1815/014'
$ws.Range("M2").Value = '17 May 2024'
$ws.Range("N2").Value = '"Date of Latest Renewal": "Not found",'
$ws.Range("O2").Value = '3 years'
$ws.Range("P2").Value = 'Store in a refrigerator (2 °C - 8 °C), do not freeze, keep the cap on the pen to protect from light'
$ws.Range("Q2").Value = 'Pre-filled pen containing 700 units of insulin icodec in 1 mL solution, 1.5 mL solution, or 3 mL solution'
$ws.Range("R2").Value = 'Prescription only'
$ws.Range("S2").Value = 'A10AE07'
$ws.Range("T2").Value = 'Novo Nordisk A/S'
$ws.Range("U2").Value = 'Glycerol, metacresol, phenol, zinc acetate, sodium chloride, hydrochloric acid, sodium hydroxide, and water for injections'
$ws.Range("V2").Value = 'awiqli-epar-product-information_en.pdf'
$ws.Range("W2").Value = '1758268326704'

# Row 3: becomes "Not found" placeholders (was row 2 content), V3 updated with leading space
$ws.Range("A3").Value = 'Not found'
$ws.Range("B3").Value = 'Not found'
$ws.Range("C3").Value = 'Not found'
$ws.Range("D3").Value = 'Not found'
$ws.Range("E3").Value = 'Not found'
$ws.Range("F3").Value = 'Not found'
$ws.Range("G3").Value = 'Not found'
$ws.Range("H3").Value = 'Not found'
$ws.Range("I3").Value = 'Not found'
$ws.Range("J3").Value = 'Not found'
$ws.Range("K3").Value = 'Not found'
$ws.Range("L3").Value = 'Not found'
$ws.Range("M3").Value = 'Not found'
$ws.Range("N3").Value = 'Not found'
$ws.Range("O3").Value = 'Not found'
$ws.Range("P3").Value = 'Not found'
$ws.Range("Q3").Value = 'Not found'
$ws.Range("R3").Value = 'Not found'
$ws.Range("S3").Value = 'Not found'
$ws.Range("T3").Value = 'Not found'
$ws.Range("U3").Value = 'Not found'
$ws.Range("V3").Value = ' https://graph.microsoft.com/v1.0/sites/slickbitai.sharepoint.com,b749c6f0-ede8-48b2-9420-ce94ca741683,876dc7c6-5b74-44d2-9d5c-a40b9e5cbf21/drive/root:/application_test1'
$ws.Range("W3").Value = '1758268326704'

# Row 4: new row with refined/expanded product data
$ws.Range("A4").Value = 'Summary of Product Characteristics, SmPC, Product Information'
$ws.Range("B4").Value = 'Awiqli'
$ws.Range("C4").Value = 'insulin icodec'
$ws.Range("D4").Value = '700 units/mL'
$ws.Range("E4").Value = 'Solution for injection in pre-filled pen'
$ws.Range("F4").Value = 'Subcutaneous use'
$ws.Range("G4").Value = 'Treatment of diabetes mellitus in adults'
$ws.Range("H4").Value = 'Once-weekly subcutaneous administration, dose adjusted based on fasting plasma glucose'
$ws.Range("I4").Value = 'Hypersensitivity to the active substance or to any of the excipients'
$ws.Range("J4").Value = 'Risk of hypoglycaemia,Switch between other insulins and insulin icodec should be done under medical supervision'
$ws.Range("K4").Value = 'Novo Nordisk A/S'
$ws.Range("L4").Value = 'EU/1/24/1815/001, EU/1/24/1815/002, EU/1/24/1815/003, EU/1/24/1815/004, EU/1/24/1815/005, EU/1/24/1815/006, EU/1/24/1815/007, EU/1/24/1815/008, EU/1/24/1815/009, EU/1/24/1815/010, EU/1/24/1815/011, EU/1/24/1815/012, EU/1/24/1815/013, EU/1/24/1815/014'
$ws.Range("M4").Value = '17 May 2024'
$ws.Range("N4").Value = '"Date of Latest Renewal": "Not found",'
$ws.Range("O4").Value = '3 years, After first opening or carried as a spare, the medicinal product may be stored for a maximum of 12 weeks'
$ws.Range("P4").Value = 'Store in a refrigerator (2°C - 8°C), Do not freeze, Keep the cap on the pen in order to protect from light'
$ws.Range("Q4").Value = '1 mL solution contains 700 units of insulin icodec, Each pre-filled pen contains 700 units of insulin icodec in 1 mL solution, Each pre-filled pen contains 1,050 units of insulin icodec in 1.5 mL solution, Each pre-filled pen contains 2,100 units of insulin icodec in 3 mL solution'
$ws.Range("R4").Value = 'Prescription only'
$ws.Range("S4").Value = 'A10AE07'
$ws.Range("T4").Value = 'Novo Nordisk A/S'
$ws.Range("U4").Value = 'Glycerol, Metacresol, Phenol, Zinc acetate, Sodium chloride, Hydrochloric acid, Sodium hydroxide, Water for injections'
$ws.Range("V4").Value = 'awiqli-epar-product-information_en.pdf'
$ws.Range("W4").Value = '1758268326704'

# Row 4 reuses the same cell style as rows 2-3 (wrap text, top-aligned)
$ws.Range("A3:W3").Copy()
$ws.Range("A4:W4").PasteSpecial(-4122)
